$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.810.30'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '2.580.97'
$ws.Range("E3").Value = '  +0.97%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = "'582.54"
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").Value = "'144.84"
$ws.Range("E6").Value = '  -1.57%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").Value = "'0.591"
$ws.Range("E8").Value = '  +1.06%  '

$ws.Range("E9").Value = '  +0.30%  '

$ws.Range("E10").Value = '  -0.38%  '

$ws.Range("E11").Value = '  -0.59%  '

$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = '  -0.27%  '

$ws.Range("D13").Value = "'27.01"
$ws.Range("E13").Value = '  -1.96%  '

$ws.Range("D14").Value = '3.041.44'
$ws.Range("E14").Value = '  +1.00%  '

$ws.Range("D15").Value = '62.687.38'
$ws.Range("E15").Value = '  -0.49%  '

$ws.Range("E16").Value = '  +0.28%  '

$ws.Range("D17").Value = '2.575.54'
$ws.Range("E17").Value = '  +0.92%  '

$ws.Range("D18").Value = "'11.23"
$ws.Range("E18").Value = '  -1.27%  '

$ws.Range("D19").Value = "'339.30"
$ws.Range("E19").Value = '  -0.33%  '

$ws.Range("D20").Value = "'4.37"
$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("E21").Value = '  -1.59%  '

$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").Value = "'5.72"
$ws.Range("E23").Value = '  +3.28%  '

$ws.Range("D24").Value = "'67.43"
$ws.Range("E24").Value = '  +2.41%  '

$ws.Range("B25").Value = 'SuiNetwork'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D25").Value = "'1.54"
$ws.Range("E25").Value = '  +3.71%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = "'0.167"
$ws.Range("E26").Value = '  -1.70%  '

$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").Value = "'1.59"
$ws.Range("E27").Value = '  -1.95%  '

$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = "'8.01"
$ws.Range("E28").Value = '  +2.35%  '

$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("D30").Value = "'8.27"
$ws.Range("E30").Value = '  -1.35%  '

$ws.Range("D31").Value = "'1.93"
$ws.Range("E31").Value = '  -2.21%  '

$ws.Range("E32").Value = '  -1.29%  '

$ws.Range("D33").Value = "'452.64"
$ws.Range("E33").Value = '  +5.77%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = "'1.63"
$ws.Range("E34").Value = '  +2.84%  '

$ws.Range("B35").Value = 'Monero'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D35").Value = "'176.28"
$ws.Range("E35").Value = '  -0.58%  '

$ws.Range("D37").Value = "'0.402"
$ws.Range("E37").Value = '  -0.91%  '

$ws.Range("D38").Value = "'18.95"
$ws.Range("E38").Value = '  -1.22%  '

$ws.Range("D39").Value = "'4.46"
$ws.Range("E39").Value = '  +0.99%  '

$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("E41").Value = '  -2.68%  '

$ws.Range("D42").Value = "'158.74"
$ws.Range("E42").Value = '  +5.04%  '

$ws.Range("E43").Value = '  -2.03%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = "'0.634"
$ws.Range("E44").Value = '  +4.88%  '

$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = "'21.21"
$ws.Range("E45").Value = '  +1.43%  '

$ws.Range("E46").Value = '  -1.52%  '

$ws.Range("D47").Value = "'0.0966"
$ws.Range("E47").Value = '  -0.78%  '

$ws.Range("D48").Value = "'0.0235"
$ws.Range("E48").Value = '  -2.37%  '

$ws.Range("D49").Value = "'18.13"
$ws.Range("E49").Value = '  -0.92%  '

$ws.Range("D50").Value = "'11.41"
$ws.Range("E50").Value = '  +0.37%  '

$ws.Range("E51").Value = '  -1.62%  '
